$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per diff
$ws.Range("B4").Value = 500
$ws.Range("B5").Value = 180

# Add new row 7: savings / 80
$ws.Range("A7").Value = "savings"
$ws.Range("B7").Value = 80
